$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'35.284.35"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.898.43"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.34%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.27%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'244.62"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.65%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.658"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +5.92%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.28%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'41.46"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -1.63%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.351"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +6.89%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'52.20"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +11.65%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0715"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +3.20%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.0996"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.58%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'2.173.65"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.43%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'12.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +4.90%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.697"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +3.31%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'1.914.17"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +3.13%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +3.19%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'35.250.77"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.46%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'71.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +2.15%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0820"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.43%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'239.30"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.60%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'12.48"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +2.01%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'4.80"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.06%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.27%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.40"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +28.81%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.31"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +2.20%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'170.74"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.80%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +6.48%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +3.81%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.95%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +4.22%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.50%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.939"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +11.79%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.20%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'4.12"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +3.22%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.64%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.02"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.09%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'1.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +3.90%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +2.17%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +4.44%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = "'16.44"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +9.64%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'0.0650"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +16.96%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'89.95"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.45%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.344.53"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.31%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'2.40"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +3.24%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'48.22"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +38.92%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'2.79"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.82%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'2.40"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.38%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'6.55"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.68%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.083.31"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +2.31%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.0696"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.28%  "
$ws.Range('E51').Style = 'Normal'
